$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 110 (the list is sorted alphabetically by
# column A, and this new dataset sorts in right above the current row 110).
$ws.Rows.Item(110).Insert()

# Fill in the new dataset's row (order matches the shared-string table
# growth in the target file: link, then name, then topics).
$ws.Range("C110").Value = "http://dx.doi.org/10.7910/DVN/FB0R8A"
$ws.Range("A110").Value = "Militant Group Electoral Participation Dataset"
$ws.Range("D110").Value = "Militant group electoral participation, rebel parties, post-conflict elections, wartime elections"
$ws.Range("B110").Value = "international relations"
$ws.Range("E110").Value = "world"
$ws.Range("F110").Value = 1970
$ws.Range("G110").Value = 2010
$ws.Range("H110").Value = "online"
$ws.Range("I110").Value = "free, no registration"

# C110 carries the dataset's link as a hyperlink, matching the rest of
# column C.
$ws.Hyperlinks.Add($ws.Range("C110"), "http://dx.doi.org/10.7910/DVN/FB0R8A")
$ws.Range("C110").Style = "Hyperlink"
